$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "EwDJUSoA"
$ws.Range("B2").Value = "26/11/2024"
$ws.Range("C2").Value = "11:00"
$ws.Range("D2").Value = "INDIA - ISL"
$ws.Range("E2").Value = "Mumbai City"
$ws.Range("F2").Value = "Punjab"
$ws.Range("G2").Value = 1.67
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.25
$ws.Range("K2").Value = 2.38
$ws.Range("L2").Value = 4.75
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 4.33
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.15
$ws.Range("S2").Value = 1.33
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.67
$ws.Range("V2").Value = 2.1
$ws.Range("W2").Value = 8.5
$ws.Range("X2").Value = 9
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 13
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 15
$ws.Range("AF2").Value = 41
$ws.Range("AG2").Value = 151
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 26
$ws.Range("AJ2").Value = 15
$ws.Range("AK2").Value = 51
$ws.Range("AL2").Value = 34
$ws.Range("AM2").Value = 34
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 8.5
$ws.Range("AP2").Value = 17
$ws.Range("AQ2").Value = 26
$ws.Range("AR2").Value = 41
$ws.Range("AS2").Value = 101
$ws.Range("AT2").Value = 3.25
$ws.Range("AU2").Value = 8
$ws.Range("AV2").Value = 51
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 23
$ws.Range("AY2").Value = 29
$ws.Range("AZ2").Value = 81
$ws.Range("BA2").Value = 81
$ws.Range("BB2").Value = 151
